$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: set "Invalid" (G3) and "Absent" (H3) to 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Rows 4-18: set "Absent" (H column) to 1
for ($r = 4; $r -le 18; $r++) {
    $ws.Cells.Item($r, 8).Value = 1
}
